# Actualización automática 2025-07-31 08:55:10
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L3").Value = 506.88
$ws1.Range("L12").Value = "1 de 10"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F3").Value = 506.88
$ws2.Range("F12").Value = 506.88
$ws2.Columns.Item(6).ColumnWidth = 11.14
